{"js": "// The document currently contains a single, empty paragraph. Insert the\n// text \"Hello\" into it (equivalent to placing the cursor in the paragraph\n// and typing \"Hello\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length > 0) {\n  const firstParagraph = paragraphs.items[0];\n  firstParagraph.insertText(\"Hello\", Word.InsertLocation.replace);\n} else {\n  body.insertParagraph(\"Hello\", Word.InsertLocation.start);\n}\n\nawait context.sync();\n", "ps1": "# The document currently contains a single, empty paragraph. Insert the\n# text \"Hello\" into it (equivalent to placing the cursor in the paragraph\n# and typing \"Hello\").\n$d = $word.ActiveDocument\n\n$para = $d.Paragraphs.First\n$para.Range.Text = \"Hello\"\n"}
